$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the italic "ब..." key-term-list paragraph entirely (it directly
#    follows the "ब" Heading2 paragraph).
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("बँधुआई, बच्चों का बलि, बतशेबा", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $d.Content.Find.Parent.Paragraphs(1)
}
$listPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "बँधुआई, बच्चों का बलि, बतशेबा*") {
        $listPara = $para
        break
    }
}
if ($listPara -ne $null) {
    $listPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Merge paragraph 6 ("This PDF version is provided under the same
#    license.") into paragraph 5, then rewrite paragraph 5's body text.
# ---------------------------------------------------------------------------
$licPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*मुख्य शब्द (Biblica) (Hindi) is based on*") {
        $licPara = $para
        break
    }
}

$mark = $d.Range($licPara.Range.End - 1, $licPara.Range.End)
$mark.Delete()

$oldSuffix = "(Hindi) is based on: Biblica Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license.This PDF version is provided under the same license."
$newSuffix = "© 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

$rSuffix = $licPara.Range.Duplicate
$fSuffix = $rSuffix.Find.Execute($oldSuffix, $false, $false, $false, $false, $false, $true, 1, $false, $newSuffix, 2)

# Rename the bold heading run itself.
$rBold = $licPara.Range.Duplicate
$fBold = $rBold.Find.Execute("मुख्य शब्द (Biblica)", $true, $false, $false, $false, $false, $true, 1, $false, "Biblica Study Notes (Key Terms)", 2)

# ---------------------------------------------------------------------------
# 3) Delete the "License Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------------
$liHeading = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "License Information*") {
        $liHeading = $para
        break
    }
}
if ($liHeading -ne $null) {
    $liHeading.Range.Delete()
}

Write-Host "done"
